$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": refresh the headline metrics now that trade #49 closed.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1197.5      # Current Capital
$summary.Range("B4").Value = -2.49       # Total P&L $
$summary.Range("B5").Value = -1.02       # Total P&L %
$summary.Range("B6").Value = 49          # Total Trades
$summary.Range("B8").Value = 26          # Losing Trades
$summary.Range("B9").Value = 38.78       # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status": MarketMaking row (row 4) mirrors the same update.
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 97.5
$status.Range("D4").Value = 49
$status.Range("E4").Value = -2.49
$status.Range("F4").Value = -2.5
$status.Range("G4").Value = 38.78

# ---------------------------------------------------------------------------
# Append the newly closed trade (#49) as row 50 on both the "All Trades"
# log and the per-strategy "MarketMaking" log - the two sheets are kept in
# sync with identical trade rows.
# ---------------------------------------------------------------------------
function Add-TradeRow {
    param($ws)

    $row = 50

    $ws.Cells.Item($row, 1).Value = 49

    # Column B holds a plain "YYYY-MM-DD" label, not a real date. A direct
    # string assignment triggers Excel's auto date-detection and turns it
    # into a date serial number with a date number format. Instead, build
    # the literal text via a throw-away text formula and paste just the
    # resulting value back in, which keeps the cell a plain, unstyled text
    # cell identical to the rest of the column.
    $ws.Range("Z1").Formula = '="2026-02-17"'
    $ws.Range("Z1").Copy()
    $ws.Cells.Item($row, 2).PasteSpecial(-4163)
    $ws.Range("Z1").Clear()

    $ws.Cells.Item($row, 3).Value = "13:28:27"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.03
    $ws.Cells.Item($row, 7).Value = 0.02
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = -33.3333
    $ws.Cells.Item($row, 10).Value = -0.01
    $ws.Cells.Item($row, 11).Value = 97.5
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.13
}

Add-TradeRow $wb.Worksheets.Item("All Trades")
Add-TradeRow $wb.Worksheets.Item("MarketMaking")
